$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet1: add the new "day" row (row 5) and move the selection
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("A5").Value = "15-Oct"

# ---------------------------------------------------------------------------
# Add the new "solution_methods" worksheet right after Sheet1
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "solution_methods"

# Data entry follows the original authoring order so new shared-string
# entries land at the same indices as the target workbook.

# A2 / A3 - Jacobian rows entered first
$ws2.Range("A2").Value = "Jacobian with loops"
$ws2.Range("A3").Value = "Jacobian vectorial"

# Header row (N / tolerance residue / time) entered next - "solver" comes later
$ws2.Range("B1").Value = "N"
$ws2.Range("C1").Value = "tolerance residue"
$ws2.Range("E1").Value = "time"
$ws2.Range("E1").Font.Bold = $true

# A4 - FiPy
$ws2.Range("A4").Value = "FiPy"

# A5 - J and F vectorial
$ws2.Range("A5").Value = "J and F vectorial "

# D1 - solver header added once the solver column was introduced
$ws2.Range("D1").Value = "solver"

# D5 - np.linalg.solve
$ws2.Range("D5").Value = "np.linalg.solve"

# A6 - J banded and F vectorial
$ws2.Range("A6").Value = "J banded and F vectorial"

# D6 - scipy.solve.banded...
$ws2.Range("D6").Value = "scipy.solve.banded with all options for performance"

# --- numeric / formatted values -------------------------------------------
$ws2.Range("B2").Value = 10
$ws2.Range("C2").Value = 0.000000001
$ws2.Range("C2").NumberFormat = "0.00E+00"
$ws2.Range("E2").Value = 1.4
$ws2.Range("E2").Font.Bold = $true

$ws2.Range("B3").Value = 10
$ws2.Range("C3").Value = 0.000000001
$ws2.Range("C3").NumberFormat = "0.00E+00"
$ws2.Range("E3").Value = 0.27
$ws2.Range("E3").Font.Bold = $true

$ws2.Range("B4").Value = 10
$ws2.Range("B4").NumberFormat = "0.00E+00"
$ws2.Range("C4").Value = 0.000000001
$ws2.Range("C4").NumberFormat = "0.00E+00"
$ws2.Range("E4").Value = 0.12
$ws2.Range("E4").Font.Bold = $true

$ws2.Range("B5").Value = 10
$ws2.Range("C5").Value = 0.000000001
$ws2.Range("C5").NumberFormat = "0.00E+00"
$ws2.Range("E5").Value = 0.1
$ws2.Range("E5").Font.Bold = $true

$ws2.Range("B6").Value = 10
$ws2.Range("C6").Value = 0.000000001
$ws2.Range("C6").NumberFormat = "0.00E+00"
$ws2.Range("E6").Value = 0.12
$ws2.Range("E6").Font.Bold = $true

# Column widths (best-fit-like)
$ws2.Columns.Item(1).ColumnWidth = 21.666666666666668
$ws2.Columns.Item(3).ColumnWidth = 15.833333333333334
$ws2.Columns.Item(4).ColumnWidth = 47.5

# ---------------------------------------------------------------------------
# Final selection state: solution_methods is the active sheet w/ H21 selected
# ---------------------------------------------------------------------------
[void]$ws2.Activate()
[void]$ws2.Range("H21").Select()

# Sheet1 selection moves to G43 (recorded while Sheet1 was active originally)
[void]$ws1.Range("G43").Select()
[void]$ws2.Activate()
